$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 598.4286
$ws.Range("I2").Value = 598.4286
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 598.4286
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -485.4286
$ws.Range("N2").ClearContents()

$ws.Range("H17").Value = 1654.1449
$ws.Range("J17").Value = 1654.1449
$ws.Range("L17").Value = 4962.4347
$ws.Range("N17").Value = -5298.4347

$ws.Range("H33").Value = 588.2308
$ws.Range("I33").Value = 618.75
$ws.Range("K33").Value = 618.75
$ws.Range("M33").Value = -389.75

$ws.Range("H61").Value = 895
$ws.Range("I61").Value = 900
$ws.Range("J61").Value = 890
$ws.Range("K61").Value = 2700
$ws.Range("L61").Value = 2670
$ws.Range("M61").Value = -2528
$ws.Range("N61").Value = -3014

$ws.Range("H74").Value = 3794.111
$ws.Range("I74").Value = 3794.111
$ws.Range("K74").Value = 3794.111
$ws.Range("M74").Value = -2858.111

$ws.Range("H77").Value = 3794.111
$ws.Range("I77").Value = 3794.111
$ws.Range("K77").Value = 18970.555
$ws.Range("M77").Value = -14290.555

$ws.Range("H88").Value = 7499
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()

$ws.Range("H91").Value = 7499
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()

$ws.Range("H116").Value = 4899.25
$ws.Range("I116").Value = 4799
$ws.Range("J116").Value = 4999.5
$ws.Range("K116").Value = 4799
$ws.Range("L116").Value = 4999.5
$ws.Range("M116").Value = -1357
$ws.Range("N116").Value = -11883.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6476.778
$ws.Range("I32").Value = 5105.1875
$ws.Range("K32").Value = 5105.1875
$ws.Range("M32").Value = -4818.1875

$ws.Range("H36").Value = 2496.25
$ws.Range("I36").Value = 661.6667
$ws.Range("K36").Value = 661.6667
$ws.Range("M36").Value = -315.6667

$ws.Range("H61").Value = 3520.2
$ws.Range("I61").Value = 3000
$ws.Range("J61").Value = 5601
$ws.Range("K61").Value = 3000
$ws.Range("L61").Value = 5601
$ws.Range("M61").Value = -2788
$ws.Range("N61").Value = -6025

$ws.Range("H122").Value = 2510.6
$ws.Range("I122").Value = 2484.842
$ws.Range("K122").Value = 7454.526
$ws.Range("M122").Value = -5004.526

$ws.Range("H131").Value = 59475
$ws.Range("J131").Value = 59475
$ws.Range("L131").Value = 59475
$ws.Range("N131").Value = -69555

$ws.Range("H136").Value = 3520.2
$ws.Range("I136").Value = 3000
$ws.Range("J136").Value = 5601
$ws.Range("K136").Value = 9000
$ws.Range("L136").Value = 16803
$ws.Range("M136").Value = -6450
$ws.Range("N136").Value = -21903

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 35000
$ws.Range("J76").Value = 35000
$ws.Range("L76").Value = 35000
$ws.Range("N76").Value = -35630

$ws.Range("H79").Value = 35000
$ws.Range("J79").Value = 35000
$ws.Range("L79").Value = 35000
$ws.Range("N79").Value = -37184

$ws.Range("H99").Value = 2000
$ws.Range("I99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("M99").Value = -502

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 6667179
$ws.Range("J22").Value = 13333701
$ws.Range("L22").Value = 13333701
$ws.Range("N22").Value = -13334401

$ws.Range("H32").Value = 3133
$ws.Range("J32").Value = 699
$ws.Range("L32").Value = 699
$ws.Range("N32").Value = -1331

$ws.Range("H59").Value = 56000
$ws.Range("I59").Value = 43333.332
$ws.Range("J59").Value = 75000
$ws.Range("K59").Value = 43333.332
$ws.Range("L59").Value = 75000
$ws.Range("M59").Value = -42188.332
$ws.Range("N59").Value = -77290

$ws.Range("H68").Value = 50000
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 50000
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H99").Value = 3000
$ws.Range("I99").Value = 3000
$ws.Range("K99").Value = 3000
$ws.Range("M99").Value = -1502

$ws.Range("H105").Value = 2712.75
$ws.Range("J105").Value = 4500
$ws.Range("L105").Value = 4500
$ws.Range("N105").Value = -7994

$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 859
$ws.Range("J11").Value = 888
$ws.Range("L11").Value = 2664
$ws.Range("N11").Value = -2944

$ws.Range("H12").Value = 345.36365
$ws.Range("J12").Value = 311.6
$ws.Range("L12").Value = 934.8000000000001
$ws.Range("N12").Value = -1280.8

$ws.Range("H113").Value = 3176
$ws.Range("J113").Value = 3093.7
$ws.Range("L113").Value = 9281.099999999999
$ws.Range("N113").Value = -13621.1

$ws.Range("H133").Value = 13304.667
$ws.Range("I133").Value = 12765.6
$ws.Range("K133").Value = 38296.8
$ws.Range("M133").Value = -33236.8

$ws.Range("H134").Value = 1127.5
$ws.Range("I134").Value = 1127.5
$ws.Range("K134").Value = 3382.5
$ws.Range("M134").Value = 1687.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6530
$ws.Range("N126").ClearContents()

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 402.69232
$ws.Range("J22").Value = 720.3333
$ws.Range("L22").Value = 720.3333
$ws.Range("N22").Value = -1310.3333

$ws.Range("H25").Value = 6000
$ws.Range("I25").Value = 6000
$ws.Range("K25").Value = 6000
$ws.Range("M25").Value = -5770

$ws.Range("H27").Value = 402.69232
$ws.Range("J27").Value = 720.3333
$ws.Range("L27").Value = 720.3333
$ws.Range("N27").Value = -934.3333

$ws.Range("H32").Value = 12999.5
$ws.Range("I32").Value = 12999.5
$ws.Range("K32").Value = 12999.5
$ws.Range("M32").Value = -12682.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 95286.5
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 95286.5
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 95286.5
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -96268.5

$ws.Range("H81").Value = 4585
$ws.Range("I81").Value = 3046.875
$ws.Range("K81").Value = 6093.75
$ws.Range("M81").Value = -5032.75

$ws.Range("H84").Value = 4585
$ws.Range("I84").Value = 3046.875
$ws.Range("K84").Value = 30468.75
$ws.Range("M84").Value = -25164.75

$ws.Range("H104").Value = 23996
$ws.Range("J104").Value = 23996
$ws.Range("L104").Value = 23996
$ws.Range("N104").Value = -30984

$ws.Range("H113").Value = 424.66666
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

$ws.Range("H130").Value = 20995
$ws.Range("J130").Value = 20995
$ws.Range("L130").Value = 20995
$ws.Range("N130").Value = -31035

Write-Host "Applied all Phantom Profits updates"